$wb = $excel.ActiveWorkbook

# --- Approver sheet: shorten the Bingo password and reorder hyperlinks ---
$approver = $wb.Worksheets.Item("Approver")

$approver.Range("B2").Value = "Bingo@12345"
$approver.Range("B3").Value = "Bingo@12345"

# Hyperlinks have to be rebuilt (this removes every hyperlink on the sheet)
# so that the B3 link now precedes the B2 link in document order, matching
# the merge result. Hyperlinks.Add's TextToDisplay also rewrites the cell
# text, so immediately restore each cell's original display text after
# recreating the klulla/kspero links.
$approver.Hyperlinks.Delete()

$approver.Hyperlinks.Add($approver.Range("A3"), "mailto:klulla@hl.com.test", "", "mailto:klulla@hl.com.test", "mailto:klulla@hl.com.test")
$approver.Range("A3").Value = "klulla@hl.com.test"

$approver.Hyperlinks.Add($approver.Range("A2"), "mailto:kspero@hl.com.test", "", "mailto:kspero@hl.com.test", "mailto:kspero@hl.com.test")
$approver.Range("A2").Value = "kspero@hl.com.test"

$approver.Hyperlinks.Add($approver.Range("B3"), "mailto:Bingo@12345")
$approver.Hyperlinks.Add($approver.Range("B2"), "mailto:Bingo@12345")

# --- Window/selection state: active tab moves from EventExp to Approver ---
$eventExp = $wb.Worksheets.Item("EventExp")
$eventExp.Range("N3").Select()

$approver.Activate()
$approver.Range("G19").Select()
